$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3600, 3, 4, 5, 6, 7, 8, 9, 11, 12, 16, 17, 20, 21, 22, 25),
    @(3601, 2, 3, 4, 6, 7, 8, 9, 11, 16, 17, 20, 21, 23, 24, 25),
    @(3602, 1, 2, 3, 4, 5, 6, 8, 9, 15, 18, 19, 20, 22, 23, 24),
    @(3603, 2, 4, 5, 9, 10, 11, 12, 14, 15, 17, 19, 20, 21, 22, 24),
    @(3604, 1, 2, 3, 6, 7, 8, 11, 13, 15, 16, 19, 21, 22, 23, 24),
    @(3605, 1, 3, 4, 7, 10, 12, 13, 14, 19, 20, 21, 22, 23, 24, 25),
    @(3606, 2, 3, 6, 7, 9, 11, 12, 14, 15, 17, 19, 20, 21, 23, 24),
    @(3607, 1, 2, 4, 5, 6, 8, 9, 13, 14, 15, 16, 18, 19, 20, 23),
    @(3608, 2, 5, 6, 8, 9, 11, 14, 16, 17, 18, 19, 20, 22, 23, 25),
    @(3609, 2, 5, 6, 8, 9, 10, 11, 12, 13, 17, 19, 20, 22, 24, 25),
    @(3610, 1, 3, 5, 7, 8, 10, 13, 14, 17, 20, 21, 22, 23, 24, 25),
    @(3611, 1, 2, 3, 4, 5, 6, 10, 11, 12, 14, 15, 17, 18, 19, 25)
)

$startRow = 530
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 0; $col -lt $rowData.Count; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $rowData[$col]
    }
}

$ws.Range("B530:P541").Select()
